$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")
$ws.Activate()

# --- Row 3 (OutcomeOfCare root concept) ---
# P3: "CarePlan  / Procedure" -> "CarePlan  / DiagnosticResult"
$ws.Range("P3").Value = "CarePlan  / DiagnosticResult"
# R3: long note -> short note (kept only the first line)
$ws.Range("R3").Value = "** OutcomeOfCare`n"
# Row 3 shrinks a lot now that the note is short
$ws.Rows.Item(3).RowHeight = 25.5

# --- Row 4 (HealthcareResult) ---
# P4: "DiagnosticReport.conclusion (HCIM Textresult)" -> new mapping text
$ws.Range("P4").Value = "CarePlan.activity:nursingIntervention.outcomeCodeableConcept Or derived profile on zib-TextResult."
# Q4: "equal" -> cleared
$ws.Range("Q4").ClearContents()

# --- Row 5 (MeasurementValue::GeneralMeasurement) ---
# Q5 text is unchanged ("-"); only the row grew taller (manual resize by author)
$ws.Rows.Item(5).RowHeight = 89.25

# --- Row 7 (Intervention::NursingIntervention) ---
# P7: "Procedure / Careplan.activity" -> new mapping text
$ws.Range("P7").Value = "Careplan.activity / DiagnosticReport.extention.partOf"
# R7: empty -> new note
$ws.Range("R7").Value = "Maybe not the most suitable extension."

# --- Selection / view bookkeeping to mirror the author's final cursor position ---
$ws.Range("R13").Select() | Out-Null
